# Applies the "cryptos" price/volume refresh described in the commit diff.
# Each row is addressed by its A1 cell refs; D-column numeric-looking text
# (e.g. "355.85") is protected with Text (@) NumberFormat so Excel keeps it
# as a literal string instead of silently re-parsing it into a Double (which
# would lose trailing zeros / introduce float noise like 355.85000000000002).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '52.270.54'
$ws.Range('E2').Value = '  -0.25%  '
# Row 3
$ws.Range('D3').Value = '2.826.56'
$ws.Range('E3').Value = '  +1.05%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '355.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.68%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.77%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.571'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.36%  '
# Row 8
$ws.Range('E8').Value = '  +0.06%  '
# Row 9
$ws.Range('E9').Value = '  +0.94%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.76%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0865'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
# Row 12
$ws.Range('E12').Value = '  +1.07%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.91%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.41%  '
# Row 15
$ws.Range('D15').Value = '3.269.69'
$ws.Range('E15').Value = '  +0.97%  '
# Row 16
$ws.Range('D16').Value = '2.827.58'
$ws.Range('E16').Value = '  +0.55%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.930'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.85%  '
# Row 18
$ws.Range('D18').Value = '52.109.96'
$ws.Range('E18').Value = '  -0.23%  '
# Row 19
$ws.Range('E19').Value = '  +4.51%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.00%  '
# Row 21
$ws.Range('E21').Value = '  +0.03%  '
# Row 22
$ws.Range('E22').Value = '  +1.42%  '
# Row 23
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '272.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.81%  '
# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
# Row 25
$ws.Range('E25').Value = '  +1.48%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.81%  '
# Row 27
$ws.Range('E27').Value = '  +0.07%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.17%  '
# Row 29
$ws.Range('E29').Value = '  +3.30%  '
# Row 30
$ws.Range('E30').Value = '  -0.82%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0491'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.08%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '52.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.18%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.84%  '
# Row 34
$ws.Range('E34').Value = '  +3.87%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.25%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0856'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.50%  '
# Row 37
$ws.Range('E37').Value = '  +0.08%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.13%  '
# Row 39
$ws.Range('E39').Value = '  -3.67%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.35%  '
# Row 41
$ws.Range('E41').Value = '  +1.77%  '
# Row 42
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '127.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.62%  '
# Row 43
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.99%  '
# Row 44
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.92%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.11%  '
# Row 46
$ws.Range('E46').Value = '  +0.41%  '
# Row 47
$ws.Range('D47').Value = '2.086.54'
$ws.Range('E47').Value = '  +0.86%  '
# Row 48
$ws.Range('E48').Value = '  -4.28%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.17%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.973'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.50%  '
# Row 51
$ws.Range('E51').Value = '  +2.47%  '
